# Auto-generated edit script to apply cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '315.46'
Set-TextValue $ws.Range('E2') '3.66%'
Set-TextValue $ws.Range('G2') '9'
Set-TextValue $ws.Range('D3') '36.07'
Set-TextValue $ws.Range('E3') '1.11%'
Set-TextValue $ws.Range('G3') '9'
Set-TextValue $ws.Range('D4') '5.174'
Set-TextValue $ws.Range('E4') '1.43%'
Set-TextValue $ws.Range('G4') '9'
Set-TextValue $ws.Range('D5') '0.08201'
Set-TextValue $ws.Range('E5') '4.46%'
Set-TextValue $ws.Range('G5') '9'
Set-TextValue $ws.Range('D6') '2.136'
Set-TextValue $ws.Range('E6') '1.47%'
Set-TextValue $ws.Range('G6') '9'
Set-TextValue $ws.Range('D7') '8.038'
Set-TextValue $ws.Range('E7') '1.45%'
Set-TextValue $ws.Range('G7') '9'
Set-TextValue $ws.Range('D8') '0.9314'
Set-TextValue $ws.Range('E8') '1.29%'
Set-TextValue $ws.Range('G8') '9'
Set-TextValue $ws.Range('D9') '0.1013'
Set-TextValue $ws.Range('E9') '3.87%'
Set-TextValue $ws.Range('G9') '9'
Set-TextValue $ws.Range('D10') '0.1876'
Set-TextValue $ws.Range('E10') '1.21%'
Set-TextValue $ws.Range('G10') '9'
Set-TextValue $ws.Range('D11') '0.09206'
Set-TextValue $ws.Range('E11') '7.53%'
Set-TextValue $ws.Range('G11') '9'
Set-TextValue $ws.Range('D12') '0.03609'
Set-TextValue $ws.Range('E12') '1.36%'
Set-TextValue $ws.Range('G12') '9'
Set-TextValue $ws.Range('D13') '0.09937'
Set-TextValue $ws.Range('E13') '-0.10%'
Set-TextValue $ws.Range('G13') '9'
Set-TextValue $ws.Range('D14') '0.001440'
Set-TextValue $ws.Range('E14') '0.11%'
Set-TextValue $ws.Range('G14') '9'
Set-TextValue $ws.Range('D15') '0.005731'
Set-TextValue $ws.Range('E15') '1.38%'
Set-TextValue $ws.Range('G15') '9'
Set-TextValue $ws.Range('D16') '3.456'
Set-TextValue $ws.Range('E16') '-0.10%'
Set-TextValue $ws.Range('G16') '9'
Set-TextValue $ws.Range('D17') '4.142'
Set-TextValue $ws.Range('E17') '0.90%'
Set-TextValue $ws.Range('G17') '9'
Set-TextValue $ws.Range('E18') '7.30%'
Set-TextValue $ws.Range('G18') '9'
Set-TextValue $ws.Range('D19') '0.3371'
Set-TextValue $ws.Range('E19') '-1.57%'
Set-TextValue $ws.Range('G19') '9'
Set-TextValue $ws.Range('E20') '1.57%'
Set-TextValue $ws.Range('G20') '9'
Set-TextValue $ws.Range('D21') '5.187'
Set-TextValue $ws.Range('E21') '-0.83%'
Set-TextValue $ws.Range('G21') '9'
Set-TextValue $ws.Range('D22') '0.2205'
Set-TextValue $ws.Range('E22') '0.08%'
Set-TextValue $ws.Range('G22') '9'
Set-TextValue $ws.Range('D23') '0.04605'
Set-TextValue $ws.Range('G23') '9'
Set-TextValue $ws.Range('D24') '0.001255'
Set-TextValue $ws.Range('E24') '1.47%'
Set-TextValue $ws.Range('G24') '9'
Set-TextValue $ws.Range('D25') '0.004717'
Set-TextValue $ws.Range('E25') '-6.72%'
Set-TextValue $ws.Range('G25') '9'
Set-TextValue $ws.Range('D26') '0.0001259'
Set-TextValue $ws.Range('E26') '-21.47%'
Set-TextValue $ws.Range('G26') '9'
Set-TextValue $ws.Range('D27') '0.0004530'
Set-TextValue $ws.Range('E27') '-4.74%'
Set-TextValue $ws.Range('G27') '9'
Set-TextValue $ws.Range('G28') '9'
Set-TextValue $ws.Range('G29') '9'
Set-TextValue $ws.Range('G30') '9'
Set-TextValue $ws.Range('G31') '9'
Set-TextValue $ws.Range('G32') '9'
Set-TextValue $ws.Range('G33') '9'
Set-TextValue $ws.Range('G34') '9'
Set-TextValue $ws.Range('G35') '9'
Set-TextValue $ws.Range('G36') '9'
Set-TextValue $ws.Range('G37') '9'
Set-TextValue $ws.Range('G38') '9'
Set-TextValue $ws.Range('D39') '0.01983'
Set-TextValue $ws.Range('E39') '8.25%'
Set-TextValue $ws.Range('G39') '9'
Set-TextValue $ws.Range('D40') '0.04910'
Set-TextValue $ws.Range('E40') '4.13%'
Set-TextValue $ws.Range('G40') '9'
Set-TextValue $ws.Range('D41') '0.007895'
Set-TextValue $ws.Range('E41') '4.41%'
Set-TextValue $ws.Range('G41') '9'
Set-TextValue $ws.Range('D42') '0.1403'
Set-TextValue $ws.Range('E42') '0.37%'
Set-TextValue $ws.Range('G42') '9'
Set-TextValue $ws.Range('D43') '0.007867'
Set-TextValue $ws.Range('E43') '1.69%'
Set-TextValue $ws.Range('G43') '9'
Set-TextValue $ws.Range('D44') '0.002125'
Set-TextValue $ws.Range('E44') '-5.25%'
Set-TextValue $ws.Range('G44') '9'
Set-TextValue $ws.Range('D45') '0.01186'
Set-TextValue $ws.Range('E45') '7.83%'
Set-TextValue $ws.Range('G45') '9'
Set-TextValue $ws.Range('D46') '0.00006560'
Set-TextValue $ws.Range('E46') '3.70%'
Set-TextValue $ws.Range('G46') '9'
Set-TextValue $ws.Range('D47') '0.00000000755'
Set-TextValue $ws.Range('E47') '0.54%'
Set-TextValue $ws.Range('G47') '9'
Set-TextValue $ws.Range('D48') '39.36'
Set-TextValue $ws.Range('E48') '-16.38%'
Set-TextValue $ws.Range('G48') '9'
Set-TextValue $ws.Range('D49') '0.001913'
Set-TextValue $ws.Range('E49') '-4.48%'
Set-TextValue $ws.Range('G49') '9'
Set-TextValue $ws.Range('D50') '0.00002114'
Set-TextValue $ws.Range('E50') '0.54%'
Set-TextValue $ws.Range('G50') '9'
Set-TextValue $ws.Range('D51') '0.0002013'
Set-TextValue $ws.Range('E51') '0.54%'
Set-TextValue $ws.Range('G51') '9'
